$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header values (C1:E1), carrying the same style as B1 ---
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)

# --- Data rows 2-21: new values for B (overwritten) plus new columns C/D/E ---
$rows = 2..21
$data = @(
    @(44.01, 26,    24,    34),
    @(50,    48,    51,    43),
    @(32,    28,    26,    22),
    @(52,    55,    56,    44),
    @(16,    16,    21,    26),
    @(36,    36,    35,    31),
    @(21,    22,    18,    21),
    @(29,    21,    23,    22),
    @(26,    46,    44,    14),
    @(24,    22,    22,    44),
    @(18,    24,    24,    49),
    @(24,    20,    26,    23),
    @(19,    21,    27,    28),
    @(8,     11,    12,    12),
    @(24,    18.1,  27,    24),
    @(56,    52,    56,    50),
    @(32.01, 22,    18,    25),
    @(8,     12.01, 11,    17),
    @(51,    54,    52.01, 49),
    @(30,    28,    35,    22.01)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $vals = $data[$i]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
}

# --- Apply the "0.000" number format (3 decimals) to the whole B:E data block ---
$ws.Range("B2:E21").NumberFormat = "0.000"
